$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values for row 9
$ws.Range("E9").Value = 92
$ws.Range("F9").Value = 0

# Update the active selection from G10 to E9
$ws.Range("E9").Select()
